# Move datasets to datasets folder
# Update PSO column (C) values for the "Training" and "Testing" worksheets
$wb = $excel.ActiveWorkbook

$wsTraining = $wb.Worksheets.Item("Training")
$trainingValues = @(
    8.7272994732618372,
    8.8713044592641381,
    8.9807986369439288,
    9.0347167337631706,
    8.5716512513415886,
    9.5480007528572859,
    8.7869830802510869,
    8.7725819126343403,
    8.8113209078371622,
    9.0037973860483405,
    8.4906799961212869,
    9.0767993752222438,
    8.7284712144319112,
    8.3948351859006376,
    9.0974352946121737,
    9.0979560942739841,
    8.4329159565046048,
    8.7689500062657544,
    8.9935463856686138,
    7.7772980565447316,
    8.2697901223039132,
    8.9071050250692902,
    8.6623127103451516,
    9.2504418853263459,
    8.8537487961303434,
    9.2121167605214591,
    8.8348827489786732,
    8.4923938877342167,
    7.8719462995498724,
    8.04433558508922
)
for ($i = 0; $i -lt $trainingValues.Length; $i++) {
    $row = 5 + $i
    $wsTraining.Cells.Item($row, 3).Value = [double]$trainingValues[$i]
}

$wsTesting = $wb.Worksheets.Item("Testing")
$testingValues = @(
    2.8046452083940321,
    2.3012554823722984,
    2.2665946875564686,
    2.0835451486535574,
    2.5590333655686379,
    1.478997732558665,
    2.31094501463518,
    2.3666078405007402,
    2.2241977574036911,
    2.3794238705716801,
    2.6550008250783965,
    1.9362067712192452,
    2.4552689442158728,
    2.7361406871608085,
    2.1302024511560829,
    2.0763356439577834,
    2.5467484154309701,
    2.2960869811848337,
    2.0149185853374294,
    3.3024855448894441,
    2.7011095309845063,
    2.4999117094730487,
    2.5911977094757055,
    1.8430918121201039,
    2.3641153681864919,
    1.9285430587232832,
    2.1981783220057078,
    2.5246228731548328,
    3.1860394748015062,
    3.0010430663772905
)
for ($i = 0; $i -lt $testingValues.Length; $i++) {
    $row = 5 + $i
    $wsTesting.Cells.Item($row, 3).Value = [double]$testingValues[$i]
}
